$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# 1. Header field changes
# ---------------------------------------------------------------
# Employee name: drop the "+B11" suffix
$ws.Range("B2").Value2 = "BAYOT ELAINE BULLON"

# A12 was stored as a text string "10/16/2014"; convert it into a real date value
$ws.Range("A12").NumberFormat = "mm/dd/yy;@"
$ws.Range("A12").Value2 = 41928

# ---------------------------------------------------------------
# 2. Fill EARNED (C) values for rows 129-135 (7 rows x 1.25)
# ---------------------------------------------------------------
$ws.Range("C129").Value2 = 1.25
$ws.Range("C130").Value2 = 1.25
$ws.Range("C131").Value2 = 1.25
$ws.Range("C132").Value2 = 1.25

# Row 133 was previously the "final templated" row (no bottom border). Its EARNED
# cell adopts the "normal" interior style (same as C132) once a real value is typed in.
$ws.Range("C132").Copy() | Out-Null
$ws.Range("C133").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C133").Value2 = 1.25

# ---------------------------------------------------------------
# 3. Extend Table1 by 14 rows (134-147) and populate the new data
# ---------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table1")
for ($i = 0; $i -lt 14; $i++) {
    $tbl.ListRows.Add() | Out-Null
}

# Copy the row-132 "interior" formatting down across 134-146, and the
# original row-133 "final row" (no bottom border) formatting onto the new
# final row, 147.
$ws.Range("A132:K132").Copy() | Out-Null
$ws.Range("A134:K146").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A133:K133").Copy() | Out-Null
$ws.Range("A147:K147").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the date number format to column A for the freshly pasted rows
$ws.Range("A134:A147").NumberFormat = "mm/dd/yy;@"

# --- row 134 : 10/1/2023, EARNED=1.25 ---
$ws.Range("A134").Value2 = 45200
$ws.Range("C134").Value2 = 1.25

# --- row 135 : 11/1/2023, EARNED=1.25 ---
$ws.Range("A135").Value2 = 45231
$ws.Range("C135").Value2 = 1.25

# --- row 136 : 12/1/2023, SP(2-0-0) taken 12/8, 11/2023 ---
$ws.Range("A136").Value2 = 45261
$ws.Range("B136").Value2 = "SP(2-0-0)"
$ws.Range("K136").Value2 = "12/8, 11/2023"

# --- row 137 : year marker "2024" ---
$ws.Range("A137").Value2 = "2024"

# --- rows 138-146 : month markers only ---
$ws.Range("A138").Value2 = 45322
$ws.Range("A139").Value2 = 45351
$ws.Range("A140").Value2 = 45382
$ws.Range("A141").Value2 = 45412
$ws.Range("A142").Value2 = 45443
$ws.Range("A143").Value2 = 45473
$ws.Range("A144").Value2 = 45504
$ws.Range("A145").Value2 = 45535
$ws.Range("A146").Value2 = 45565

# --- row 147 : final templated row ---
$ws.Range("A147").Value2 = 45596

# ---------------------------------------------------------------
# 4. Make sure the EARNED " (G) formulas exist for every new row and
#    recalc the whole workbook so cached formula results match.
# ---------------------------------------------------------------
$ws.Range("G134").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G135").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G136").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G137").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G138").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G139").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G140").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G141").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G142").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G143").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G144").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G145").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G146").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G147").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$excel.CalculateFull()
